$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert new row 6: Grand Island (JBS) plant ---
$ws.Rows.Item(6).Insert()
$ws.Range("R6").Clear()
$ws.Range("A6").Value = "Grand Island"
$ws.Range("B6").Value = "JBS"
$ws.Range("C6").Value = "USA"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "Beef"
$ws.Range("G6").Value = "Plant"
$ws.Range("S6").Value = "555 South Stuhr Road, Grand Island, NE 68801"
$ws.Range("T6").Value = 68801

# --- Insert new row 21: Willmar (Jennie-O) plant ---
$ws.Rows.Item(21).Insert()
$ws.Range("A21").Value = "Willmar"
$ws.Range("B21").Value = "Jennie-O"
$ws.Range("C21").Value = "USA"
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = "Beef"
$ws.Range("G21").Value = "Plant"
$ws.Range("S21").Value = "2505 Willmar Ave SW, Willmar, MN 56201"
$ws.Range("T21").Value = 56201

# --- Fix up U1:V1 header style (drop the stray fill-applied style) ---
$ws.Range("A1").Copy()
$ws.Range("U1:V1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update selection / scroll position ---
$ws.Range("A21").Select()
